$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "DONE"
$ws.Range("A15").Value = "Done"
$ws.Range("A18").Value = "Done"
$ws.Range("A21").Value = "Done"

$ws.Range("A15").Select()
